$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a weekly price log for "Naranja" (orange) at
# "Vega Monumental Concepción". Two new rows of data (Navel Late,
# Primera / Segunda, dated 2021-12-21 / serial 44551) were added right
# after the header block for this market, pushing every following
# record down by two rows (old row N -> new row N+2).
#
# Strategy: insert two blank rows at row 164 (Excel shifts rows 164..195
# down to 166..197, which already reproduces the "shift by two" pattern
# seen in the diff), then seed rows 164/165 with a duplicate of the
# (now shifted) row 166/167 so every column/style inherits correctly,
# and finally overwrite the handful of columns that actually differ for
# the two brand-new records.

$ws.Rows.Item(164).Insert()
$ws.Rows.Item(164).Insert()

for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item(164, $c).Value2 = $ws.Cells.Item(166, $c).Value2
    $ws.Cells.Item(165, $c).Value2 = $ws.Cells.Item(167, $c).Value2
}

# Match the date cell's number format (column D uses a date/time style).
$ws.Range("D164:D165").NumberFormat = $ws.Range("D166").NumberFormat

# Row 164: Navel Late / Primera
$ws.Cells.Item(164, 4).Value2 = 44551   # Fecha
$ws.Cells.Item(164, 11).Value2 = "Navel Late"   # Variedad
$ws.Cells.Item(164, 12).Value2 = "Primera"      # Calidad
$ws.Cells.Item(164, 13).Value2 = 450    # Volumen
$ws.Cells.Item(164, 14).Value2 = 8500   # Precio minimo
$ws.Cells.Item(164, 15).Value2 = 9000   # Precio maximo
$ws.Cells.Item(164, 16).Value2 = 8778   # Precio promedio ponderado
$ws.Cells.Item(164, 17).Value2 = "$/caja 15 kilos empedrada"  # Unidad de comercializacion
$ws.Cells.Item(164, 18).Value2 = "Región de O'Higgins"        # Origen
$ws.Cells.Item(164, 19).Value2 = 585    # Precio $/Kg
$ws.Cells.Item(164, 20).Value2 = 15     # Kg / unidad

# Row 165: Navel Late / Segunda
$ws.Cells.Item(165, 4).Value2 = 44551
$ws.Cells.Item(165, 11).Value2 = "Navel Late"
$ws.Cells.Item(165, 12).Value2 = "Segunda"
$ws.Cells.Item(165, 13).Value2 = 350
$ws.Cells.Item(165, 14).Value2 = 7500
$ws.Cells.Item(165, 15).Value2 = 8000
$ws.Cells.Item(165, 16).Value2 = 7714
$ws.Cells.Item(165, 17).Value2 = "$/caja 15 kilos empedrada"
$ws.Cells.Item(165, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(165, 19).Value2 = 514
$ws.Cells.Item(165, 20).Value2 = 15
